# NewsWebsiteData.xlsx update:
#  - "A specific story and author" update: fill in the (previously blank) numeric
#    id column (A) on the "journalists" and "types" sheets with a simple 1..7
#    row counter ("added icons" / numbering next to each row).
#  - Switch the active sheet/selection to reflect where the author ended up
#    working: "articles" sheet, cell H11 selected (was "journalists"!D16).
#  - Update the remembered selection on "journalists" (-> B11) and
#    "types" (-> C4) for when the user returns to them.

$wb = $excel.ActiveWorkbook

$articles    = $wb.Worksheets.Item("articles")
$journalists = $wb.Worksheets.Item("journalists")
$types       = $wb.Worksheets.Item("types")

# --- Fill column A (id numbers) for journalists!A2:A8 and types!A2:A8 ---
for ($row = 2; $row -le 8; $row++) {
    $journalists.Cells.Item($row, 1).Value = $row - 1
}

for ($row = 2; $row -le 8; $row++) {
    $types.Cells.Item($row, 1).Value = $row - 1
}

# --- Update remembered selections on the non-active sheets ---
$journalists.Range("B11").Select()
$types.Range("C4").Select()

# --- Make "articles" the active sheet with H11 selected ---
$articles.Activate()
$articles.Range("H11").Select()
